# Refined metadata to be additional tab
#
# 1) Add a new "metadata" worksheet after the existing "data" sheet.
# 2) Populate it with one metadata record describing the panel query.
# 3) Re-stamp the "time_taken" column on the "data" sheet with the new
#    query timestamp (panel_query_time changed, so every row's
#    time_taken on data moves forward by the same delta as F2 on metadata).

$wb = $excel.ActiveWorkbook
$ds = $wb.Worksheets.Item("data")

# --- add the metadata sheet right after "data" --------------------------
$ws = $wb.Worksheets.Add($null, $ds)
$ws.Name = "metadata"

# --- header row -----------------------------------------------------------
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

$hdr = $ws.Range("B1:G1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160

# --- data row ---------------------------------------------------------
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Brain channelopathy"
$ws.Range("C2").Value = 90
$ws.Range("D2").Value = "'1.60"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "2021-06-21T09:41:40.073037Z"
$ws.Range("F2").Value = "2021-10-05 14:19:20.938770"
$ws.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/90/?format=json"

$a2 = $ws.Range("A2")
$a2.Font.Bold = $true
$a2.Borders.LineStyle = 1
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160

# --- re-stamp the data sheet's time_taken column --------------------------
$newTimes = @(
    "2021-10-05 14:19:20.942614",
    "2021-10-05 14:19:20.942622",
    "2021-10-05 14:19:20.942625",
    "2021-10-05 14:19:20.942628",
    "2021-10-05 14:19:20.942631",
    "2021-10-05 14:19:20.942634",
    "2021-10-05 14:19:20.942636",
    "2021-10-05 14:19:20.942639",
    "2021-10-05 14:19:20.942641",
    "2021-10-05 14:19:20.942644",
    "2021-10-05 14:19:20.942647",
    "2021-10-05 14:19:20.942649",
    "2021-10-05 14:19:20.942651",
    "2021-10-05 14:19:20.942654",
    "2021-10-05 14:19:20.942656",
    "2021-10-05 14:19:20.942659",
    "2021-10-05 14:19:20.942662",
    "2021-10-05 14:19:20.942664",
    "2021-10-05 14:19:20.942667",
    "2021-10-05 14:19:20.942669",
    "2021-10-05 14:19:20.942672",
    "2021-10-05 14:19:20.942674",
    "2021-10-05 14:19:20.942677",
    "2021-10-05 14:19:20.942679",
    "2021-10-05 14:19:20.942682",
    "2021-10-05 14:19:20.942685"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = 2 + $i
    $ds.Range("F$row").Value = $newTimes[$i]
}

# --- keep "data" as the active sheet (bookViews activeTab stays 0) --------
$ds.Activate()
